$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the formatting (bold / border / centered) already used by the other header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the I and J columns for the data rows.
# Column I is a constant 1, column J duplicates whatever value is already in column H.
for ($r = 2; $r -le 37; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Text
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
